$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current row 379 (old data shifts down by 2,
# row count grows from 432 to 434, matching the diff's dimension change).
$ws.Rows("379:380").Insert()

# Row 379: new weekly data point (Primera quality)
$ws.Cells.Item(379, 1).Value  = 9
$ws.Cells.Item(379, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(379, 3).Value  = "Metropolitana"
$ws.Cells.Item(379, 4).Value  = 45034
$ws.Cells.Item(379, 5).Value  = 13
$ws.Cells.Item(379, 6).Value  = 100112017
$ws.Cells.Item(379, 7).Value  = "Apio"
$ws.Cells.Item(379, 8).Value  = "Americana (o)"
$ws.Cells.Item(379, 9).Value  = "Primera"
$ws.Cells.Item(379, 10).Value = 70
$ws.Cells.Item(379, 11).Value = 8000
$ws.Cells.Item(379, 12).Value = 9000
$ws.Cells.Item(379, 13).Value = 8500
$ws.Cells.Item(379, 14).Value = "`$/docena de matas"
$ws.Cells.Item(379, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(379, 16).Value = 1417
$ws.Cells.Item(379, 17).Value = 6
$ws.Cells.Item(379, 18).Value = "Hortaliza"

# Row 380: new weekly data point (Segunda quality)
$ws.Cells.Item(380, 1).Value  = 9
$ws.Cells.Item(380, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(380, 3).Value  = "Metropolitana"
$ws.Cells.Item(380, 4).Value  = 45034
$ws.Cells.Item(380, 5).Value  = 13
$ws.Cells.Item(380, 6).Value  = 100112017
$ws.Cells.Item(380, 7).Value  = "Apio"
$ws.Cells.Item(380, 8).Value  = "Americana (o)"
$ws.Cells.Item(380, 9).Value  = "Segunda"
$ws.Cells.Item(380, 10).Value = 52
$ws.Cells.Item(380, 11).Value = 7000
$ws.Cells.Item(380, 12).Value = 7000
$ws.Cells.Item(380, 13).Value = 7000
$ws.Cells.Item(380, 14).Value = "`$/docena de matas"
$ws.Cells.Item(380, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(380, 16).Value = 1167
$ws.Cells.Item(380, 17).Value = 6
$ws.Cells.Item(380, 18).Value = "Hortaliza"
